$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Capture the width of the existing "In Advance" column (M) so the newly
# inserted column can be given a matching custom width.
$srcColumnWidth = $ws.Range("M1").ColumnWidth

# Insert a new blank column before the "Late" column (N), shifting the
# "Late" and "Outstanding" columns one place to the right (to O and Q,
# the latter leaving a blank gap column in between, matching the
# pre-existing layout where an unused column separated the data).
$ws.Columns("N").Insert()

# Give the freshly inserted column a custom width similar to its neighbor.
$ws.Range("N1").ColumnWidth = $srcColumnWidth

# Update the active selection on the sheet, as left by the editor.
$ws.Range("S8").Select()
